$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.005.15"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "'1.873.37"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'305.63"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.5088"
$ws.Range("E7").Value = "  -0.51%  "

$ws.Range("D8").Value = "'0.3671"
$ws.Range("E8").Value = "  -1.99%  "

$ws.Range("D9").Value = "'0.07215"
$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("D10").Value = "'0.8961"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.872.72"
$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07528"
$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'95.38"
$ws.Range("E14").Value = "  +6.60%  "

$ws.Range("D15").Value = "'5.252"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").Value = "'14.25"
$ws.Range("E18").Value = "  +0.97%  "

$ws.Range("D19").Value = "'0.9994"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").Value = "'27.033.45"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").Value = "'5.030"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "'2.101.77"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("E23").Value = "  -1.23%  "

$ws.Range("D24").Value = "'6.403"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").Value = "'148.37"
$ws.Range("E25").Value = "  +0.34%  "

$ws.Range("E26").Value = "  -3.04%  "

$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("E28").Value = "  -1.14%  "

$ws.Range("D29").Value = "'113.42"
$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("D30").Value = "'4.726"

$ws.Range("D31").Value = "'4.739"
$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("D32").Value = "'0.09154"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").Value = "'0.05120"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").Value = "'0.7516"
$ws.Range("E34").Value = "  +3.39%  "

$ws.Range("D35").Value = "'2.975"
$ws.Range("E35").Value = "  -2.52%  "

$ws.Range("D36").Value = "'1.163"
$ws.Range("E36").Value = "  +0.68%  "

$ws.Range("D37").Value = "'3.238"
$ws.Range("E37").Value = "  +6.40%  "

$ws.Range("D38").Value = "'2.556"
$ws.Range("E38").Value = "  +2.48%  "

$ws.Range("D39").Value = "'0.5659"
$ws.Range("E39").Value = "  +6.09%  "

$ws.Range("D40").Value = "'0.02002"
$ws.Range("E40").Value = "  -1.85%  "

$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").Value = "'6.644"
$ws.Range("E42").Value = "  +1.12%  "

$ws.Range("D43").Value = "'115.67"
$ws.Range("E43").Value = "  -1.13%  "

$ws.Range("D44").Value = "'8.587"
$ws.Range("E44").Value = "  +3.84%  "

$ws.Range("D45").Value = "'0.1477"
$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("D46").Value = "'0.4785"
$ws.Range("E46").Value = "  +3.23%  "

$ws.Range("D47").Value = "'0.9996"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").Value = "'10.11"
$ws.Range("E48").Value = "  +1.14%  "

$ws.Range("D49").Value = "'1.573"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").Value = "'36.96"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").Value = "'63.24"
$ws.Range("E51").Value = "  -1.07%  "

